# Generate Report for Handoff
# - Flip status strings from "Handed back: in sync with en-US" to "Ready for handoff"
# - Refresh the handoff timestamps
# - Narrow the now-shorter status columns

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Status text: "Handed back: in sync with en-US" -> "Ready for handoff"
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# Refreshed handoff timestamps
$overview.Range("G2").Value = "2016-08-24 21:01:59"
$zhcn.Range("H2").Value = "2016-08-24 21:01:53"
$dede.Range("H2").Value = "2016-08-24 21:01:59"

# Columns shrink now that the status text is shorter
# (Excel.ColumnWidth is character-unit; the closest value that lands the
# underlying sheet column width next to the target ~17.216 is ~16.333)
$overview.Range("E1").ColumnWidth = 16.3333333333333
$overview.Range("F1").ColumnWidth = 16.3333333333333
$zhcn.Range("C1").ColumnWidth = 16.3333333333333
$dede.Range("C1").ColumnWidth = 16.3333333333333
